$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet: "No_ANT" -> "NO_ANT"
$ws.Name = "PickAndPlace_PCB V1.1_NO_ANT_20"

# Row 11 "Comment" column changes from "pi" to "Pi Header"
$ws.Cells.Item(11, 14).Value = "Pi Header"

# Add new "Supplier Part" column (O / column 15)
$ws.Cells.Item(1, 15).Value = "Supplier Part"

$ws.Cells.Item(2, 15).Value  = "C124378"
$ws.Cells.Item(3, 15).Value  = "C124378"
$ws.Cells.Item(4, 15).Value  = "C96446"
$ws.Cells.Item(5, 15).Value  = "C96446"
$ws.Cells.Item(6, 15).Value  = "C14663"
$ws.Cells.Item(7, 15).Value  = "C14663"
$ws.Cells.Item(8, 15).Value  = "C5120765"
$ws.Cells.Item(9, 15).Value  = "C28323"
$ws.Cells.Item(10, 15).Value = "C16133"
$ws.Cells.Item(11, 15).Value = "C2977589"
$ws.Cells.Item(12, 15).Value = "C411294"
$ws.Cells.Item(13, 15).Value = "C11702"
$ws.Cells.Item(14, 15).Value = "C11702"
$ws.Cells.Item(15, 15).Value = "C86038"
$ws.Cells.Item(16, 15).Value = "C86038"
# Rows 17 and 18 have no supplier part (left blank)

# Match column width of new column O to the rest (A:N are width 20)
$ws.Range("O1:O18").ColumnWidth = 19.16666666666667
